$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 404.625
$ws.Range("I2").Value = 404.625
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 404.625
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -291.625

$ws.Range("H21").Value = 298.5
$ws.Range("I21").Value = 298.5
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 298.5
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 169.5

$ws.Range("H23").Value = 298.5
$ws.Range("I23").Value = 298.5
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 298.5
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -64.5

$ws.Range("H33").Value = 244.05556
$ws.Range("I33").Value = 199.6875
$ws.Range("J33").Value = 599
$ws.Range("K33").Value = 199.6875
$ws.Range("L33").Value = 599
$ws.Range("M33").Value = 29.3125
$ws.Range("N33").Value = -1057

$ws.Range("H38").Value = 614.4
$ws.Range("I38").Value = 614.4
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 1843.2
$ws.Range("L38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("N38").Value = -1471.2

$ws.Range("H41").Value = 647.7143
$ws.Range("I41").Value = 647.7143
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 647.7143
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -207.7143

$ws.Range("H51").Value = 7824.875
$ws.Range("I51").Value = 2650
$ws.Range("J51").Value = 12999.75
$ws.Range("K51").Value = 2650
$ws.Range("L51").Value = 12999.75
$ws.Range("M51").Value = -2166
$ws.Range("N51").Value = -13967.75

$ws.Range("H70").Value = 1666
$ws.Range("I70").Value = 1999
$ws.Range("J70").Value = 1499.5
$ws.Range("K70").Value = 5997
$ws.Range("L70").Value = 4498.5
$ws.Range("M70").Value = -5727
$ws.Range("N70").Value = -5038.5

$ws.Range("H73").Value = 1666
$ws.Range("I73").Value = 1999
$ws.Range("J73").Value = 1499.5
$ws.Range("K73").Value = 5997
$ws.Range("L73").Value = 4498.5
$ws.Range("M73").Value = -5061
$ws.Range("N73").Value = -6370.5

$ws.Range("H86").Value = 2364.25
$ws.Range("I86").Value = 2570.7144
$ws.Range("J86").Value = 2075.2
$ws.Range("K86").Value = 2570.7144
$ws.Range("L86").Value = 2075.2
$ws.Range("M86").Value = -1447.7144
$ws.Range("N86").Value = -4321.2

$ws.Range("H89").Value = 2364.25
$ws.Range("I89").Value = 2570.7144
$ws.Range("J89").Value = 2075.2
$ws.Range("K89").Value = 12853.572
$ws.Range("L89").Value = 10376
$ws.Range("M89").Value = -7237.572
$ws.Range("N89").Value = -21608

$ws.Range("H132").Value = 2664.8823
$ws.Range("I132").Value = 2581.4375
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 7744.3125
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -5214.3125
$ws.Range("N132").Value = -17060

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 3779.8333
$ws.Range("I19").Value = 2194.75
$ws.Range("J19").Value = 6950
$ws.Range("K19").Value = 2194.75
$ws.Range("L19").Value = 6950
$ws.Range("M19").Value = -1965.75
$ws.Range("N19").Value = -7408

$ws.Range("H45").Value = 945.8570999999999
$ws.Range("I45").Value = 821.4
$ws.Range("J45").Value = 1257
$ws.Range("K45").Value = 821.4
$ws.Range("L45").Value = 1257
$ws.Range("M45").Value = -444.4
$ws.Range("N45").Value = -2011

$ws.Range("H74").Value = 1749.5454
$ws.Range("I74").Value = 1666.2222
$ws.Range("J74").Value = 2124.5
$ws.Range("K74").Value = 1666.2222
$ws.Range("L74").Value = 2124.5
$ws.Range("M74").Value = -792.2221999999999
$ws.Range("N74").Value = -3872.5

$ws.Range("H77").Value = 1749.5454
$ws.Range("I77").Value = 1666.2222
$ws.Range("J77").Value = 2124.5
$ws.Range("K77").Value = 8331.110999999999
$ws.Range("L77").Value = 10622.5
$ws.Range("M77").Value = -3963.110999999999
$ws.Range("N77").Value = -19358.5

$ws.Range("H131").Value = 50000
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 50000
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 50000
$ws.Range("N131").Value = -60080

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2558.0833
$ws.Range("I20").Value = 1629.8572
$ws.Range("J20").Value = 3857.6
$ws.Range("K20").Value = 1629.8572
$ws.Range("L20").Value = 3857.6
$ws.Range("M20").Value = -1382.8572
$ws.Range("N20").Value = -4351.6

$ws.Range("H105").Value = 71432344
$ws.Range("I105").Value = 2423.8
$ws.Range("J105").Value = 111115630
$ws.Range("K105").Value = 2423.8
$ws.Range("L105").Value = 111115630
$ws.Range("M105").Value = -676.8000000000002
$ws.Range("N105").Value = -111119124

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 80000
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 80000
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 80000
$ws.Range("N20").Value = -80472

$ws.Range("H30").Value = 80000
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 80000
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 80000
$ws.Range("N30").Value = -80182

$ws.Range("H31").Value = 2184.5
$ws.Range("I31").Value = 1796.5714
$ws.Range("J31").Value = 4900
$ws.Range("K31").Value = 1796.5714
$ws.Range("L31").Value = 4900
$ws.Range("M31").Value = -1501.5714
$ws.Range("N31").Value = -5490

$ws.Range("H34").Value = 2184.5
$ws.Range("I34").Value = 1796.5714
$ws.Range("J34").Value = 4900
$ws.Range("K34").Value = 1796.5714
$ws.Range("L34").Value = 4900
$ws.Range("M34").Value = -1594.5714
$ws.Range("N34").Value = -5304

$ws.Range("H58").Value = 3603.5
$ws.Range("I58").Value = 2804.8333
$ws.Range("J58").Value = 5999.5
$ws.Range("K58").Value = 2804.8333
$ws.Range("L58").Value = 5999.5
$ws.Range("M58").Value = -2601.8333
$ws.Range("N58").Value = -6405.5

$ws.Range("H128").Value = 80000
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 80000
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 80000
$ws.Range("N128").Value = -89960

$ws.Range("H136").Value = 3603.5
$ws.Range("I136").Value = 2804.8333
$ws.Range("J136").Value = 5999.5
$ws.Range("K136").Value = 8414.499899999999
$ws.Range("L136").Value = 17998.5
$ws.Range("M136").Value = -5864.499899999999
$ws.Range("N136").Value = -23098.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 4383.2856
$ws.Range("I34").Value = 394.66666
$ws.Range("J34").Value = 7374.75
$ws.Range("K34").Value = 1183.99998
$ws.Range("L34").Value = 22124.25
$ws.Range("M34").Value = -1099.99998
$ws.Range("N34").Value = -22292.25

$ws.Range("H70").Value = 11509.7
$ws.Range("I70").Value = 2724.25
$ws.Range("J70").Value = 17366.666
$ws.Range("K70").Value = 8172.75
$ws.Range("L70").Value = 52099.99800000001
$ws.Range("M70").Value = -7857.75
$ws.Range("N70").Value = -52729.99800000001

$ws.Range("H73").Value = 11509.7
$ws.Range("I73").Value = 2724.25
$ws.Range("J73").Value = 17366.666
$ws.Range("K73").Value = 8172.75
$ws.Range("L73").Value = 52099.99800000001
$ws.Range("M73").Value = -7080.75
$ws.Range("N73").Value = -54283.99800000001

$ws.Range("H112").Value = 10211.3
$ws.Range("I112").Value = 1056.5
$ws.Range("J112").Value = 12500
$ws.Range("K112").Value = 3169.5
$ws.Range("L112").Value = 37500
$ws.Range("M112").Value = -2061.5
$ws.Range("N112").Value = -39716

$ws.Range("H113").Value = 1114.6428
$ws.Range("I113").Value = 369.8
$ws.Range("J113").Value = 1528.4445
$ws.Range("K113").Value = 1109.4
$ws.Range("L113").Value = 4585.333500000001
$ws.Range("M113").Value = 1060.6
$ws.Range("N113").Value = -8925.333500000001

$ws.Range("H137").Value = 5387.5386
$ws.Range("I137").Value = 4425
$ws.Range("J137").Value = 5562.5454
$ws.Range("K137").Value = 13275
$ws.Range("L137").Value = 16687.6362
$ws.Range("M137").Value = -8175
$ws.Range("N137").Value = -26887.6362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8365.333000000001
$ws.Range("I70").Value = 7961.6665
$ws.Range("J70").Value = 8499.888999999999
$ws.Range("K70").Value = 7961.6665
$ws.Range("L70").Value = 8499.888999999999
$ws.Range("M70").Value = -7691.6665
$ws.Range("N70").Value = -9039.888999999999

$ws.Range("H73").Value = 8365.333000000001
$ws.Range("I73").Value = 7961.6665
$ws.Range("J73").Value = 8499.888999999999
$ws.Range("K73").Value = 7961.6665
$ws.Range("L73").Value = 8499.888999999999
$ws.Range("M73").Value = -7025.6665
$ws.Range("N73").Value = -10371.889

$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()

$ws.Range("H102").Value = 2364.2856
$ws.Range("I102").Value = 2364.2856
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2364.2856
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -742.2856000000002

$ws.Range("H126").Value = 2179.1
$ws.Range("I126").Value = 1968.8334
$ws.Range("J126").Value = 2494.5
$ws.Range("K126").Value = 5906.5002
$ws.Range("L126").Value = 7483.5
$ws.Range("M126").Value = -3436.5002
$ws.Range("N126").Value = -12423.5

$ws.Range("H128").Value = 20890
$ws.Range("I128").Value = 20890
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 20890
$ws.Range("L128").Value = 0
$ws.Range("M128").Value = -15910
$ws.Range("N128").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1008.8
$ws.Range("I93").Value = 961
$ws.Range("J93").Value = 1200
$ws.Range("K93").Value = 961
$ws.Range("L93").Value = 1200
$ws.Range("M93").Value = 287
$ws.Range("N93").Value = -3696

$ws.Range("H128").Value = 80000
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 80000
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 80000
$ws.Range("N128").Value = -89960

$ws.Range("H136").Value = 8500
$ws.Range("I136").Value = 1000
$ws.Range("J136").Value = 12250
$ws.Range("K136").Value = 3000
$ws.Range("L136").Value = 36750
$ws.Range("M136").Value = -450
$ws.Range("N136").Value = -41850

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 5000
$ws.Range("I20").Value = 5000
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 5000
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -4760

$ws.Range("H130").Value = 36333.332
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 36333.332
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 36333.332
$ws.Range("N130").Value = -46373.332
